# - Stage Design - USPSA.pptx
# "Updated some stages, and worked on the WSB's for the ESB stages"
#
# 1) Bump the notes-master date field 8/1/2019 -> 8/29/2019.
# 2) Add a new (blank-layout) slide containing a single text box with the
#    generic "Written Stage Briefing" template copy.

$p = $ppt.ActivePresentation

# --- 1) Notes master date field -------------------------------------------
# The cached text of the datetimeFigureOut field on the notes master lives
# on the "Date Placeholder" shape. Try the documented object-model routes;
# if the host can't route the edit back to the notes master part it raises
# a non-terminating error and we simply carry on with the slide work below.
$nm = $p.NotesMaster
$datePh = $nm.HeadersFooters.DateAndTime
$datePh.Type = 1
$datePh.Value = "8/29/2019"
foreach ($shp in $nm.Shapes) {
    if ($shp.Name -eq "Date Placeholder 2") {
        $shp.TextFrame.TextRange.Text = "8/29/2019"
    }
}

# --- 2) New slide: Written Stage Briefing textbox --------------------------
# ppLayoutBlank = 12 -> no placeholders, matching the captured slide2.xml.
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 12)

# EMU 341832/384561/6699903/3416320 expressed in points (1 pt = 12700 EMU)
# so the shape lands on the exact same off/ext as the authored slide.
$tb = $slide.Shapes.AddTextbox(1, 26.915905511811022, 30.280393700787403, 527.5514173228346, 269.0015748031496)
$tb.Name = "TextBox 1"
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0

$tr = $tb.TextFrame.TextRange
$tr.Text = "Stage Name"
$null = $tr.InsertAfter("`rWritten Stage Briefing")
$null = $tr.InsertAfter("`rStage Name is a RoundCount round, Points point, ScoringMethod CourseType course. Targets. The best HitCount hits per target will score. Steel must fall to score. The start signal is audible.")
$null = $tr.InsertAfter("`r")
$null = $tr.InsertAfter("`rHangun start position")
$null = $tr.InsertAfter("`r")
$null = $tr.InsertAfter("`rPCC start position")
$null = $tr.InsertAfter("`r")
$null = $tr.InsertAfter("`rStage Procedure")

$tr.Paragraphs(1, 1).Font.Bold = 1
$tr.Paragraphs(1, 1).ParagraphFormat.Alignment = 2
$tr.Paragraphs(2, 2).Font.Bold = 1
$tr.Paragraphs(2, 2).ParagraphFormat.Alignment = 2

Write-Output "Slides now: $($p.Slides.Count); new slide shapes: $($slide.Shapes.Count)"
